$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "Louisiana Plus", 97, 0.78, 0.76, 0, "04-07")
    ,@(3, "Louisiana Plus", 337, 0.72, 0.71, 1, "12-03")
    ,@(4, "East Canada", 177, 0.84, 0.82, 0, "06-26")
    ,@(5, "East Canada", 257, 0.86, 0.86, 1, "09-14")
    ,@(6, "Northeast Plus", 157, 0.79, 0.76, 0, "06-06")
    ,@(7, "Northeast Plus", 277, 0.79, 0.8, 1, "10-04")
    ,@(8, "Missouri Plus", 147, 0.76, 0.76, 0, "05-27")
    ,@(9, "Missouri Plus", 282, 0.78, 0.74, 1, "10-09")
    ,@(10, "Florida State", 37, 0.75, 0.73, 0, "02-06")
    ,@(11, "Florida State", 362, 0.75, 0.72, 1, "12-28")
    ,@(12, "Alaska State", 212, 0.84, 0.85, 0, "07-31")
    ,@(13, "Alaska State", 217, 0.86, 0.85, 1, "08-05")
    ,@(14, "Near DC", 137, 0.77, 0.78, 0, "05-17")
    ,@(15, "Near DC", 297, 0.8, 0.82, 1, "10-24")
    ,@(16, "Minnesota Plus", 177, 0.85, 0.78, 0, "06-26")
    ,@(17, "Minnesota Plus", 252, 0.82, 0.82, 1, "09-09")
    ,@(18, "Georgia Plus", 117, 0.83, 0.79, 0, "04-27")
    ,@(19, "Georgia Plus", 302, 0.79, 0.79, 1, "10-29")
    ,@(20, "Indiana Plus", 157, 0.79, 0.76, 0, "06-06")
    ,@(21, "Indiana Plus", 257, 0.84, 0.79, 1, "09-14")
    ,@(22, "New Mexico Plus", 112, 0.66, 0.65, 0, "04-22")
    ,@(23, "New Mexico Plus", 297, 0.63, 0.64, 1, "10-24")
    ,@(24, "Oregon Plus", 177, 0.78, 0.8100000000000001, 0, "06-26")
    ,@(25, "Oregon Plus", 267, 0.82, 0.85, 1, "09-24")
    ,@(26, "California Plus", 147, 0.77, 0.79, 0, "05-27")
    ,@(27, "California Plus", 311, 0.76, 0.76, 1, "11-07")
    ,@(28, "California Plus", 312, 0.76, 0.76, 1, "11-08")
    ,@(29, "Hawaii State", 72, 0.8100000000000001, 0.83, 0, "03-13")
    ,@(30, "Hawaii State", 365, 0.67, 0.67, 1, "12-31")
    ,@(31, "Wyoming Plus", 157, 0.71, 0.6899999999999999, 0, "06-06")
    ,@(32, "Wyoming Plus", 272, 0.67, 0.67, 1, "09-29")
    ,@(33, "Puerto Rico", 32, 0.14, 0.2, 0, "02-01")
    ,@(34, "Puerto Rico", 365, 0.13, 0.1, 1, "12-31")
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
    $ws.Cells.Item($r, 5).Value = $item[5]
    $ws.Cells.Item($r, 7).Value = $item[6]
}